$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix shared string text: "United_States" -> "United States"
$ws.Range("C10:C17").Value = "United States"

# Fix shared string text: "England&Wales" -> "England & Wales"
$ws.Range("C18:C25").Value = "England & Wales"

# Update selection to C9
$ws.Range("C9").Select()
